# sửa login.html nhưng lỗi đăng nhập
# Append an empty paragraph plus two new paragraphs with SQL text after
# the existing "   Và jdk 23" paragraph (the last paragraph in the doc).

$d = $word.ActiveDocument

# 1) Empty paragraph right after "   Và jdk 23"
$d.Paragraphs.Last.Range.InsertParagraphAfter()

# 2) "INSERT INTO ..." paragraph
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.InsertBefore("INSERT INTO nguoidungs (username, password, role, avatar, email)")

# 3) "VALUES (...)" paragraph
$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.InsertBefore("VALUES ('giaovu1', '123456', 'giaovu', 'avatar.png', 'giaovu1@example.com');")
